# Sprint Log.xlsx update
# - Bumps the Sprint 1 total points (F2) from 8 to 17 (burndown chart cache follows).
# - Marks the "Encrypt User Passwords" / "Implement Login Functionality" /
#   "Add Login UI" rows (10-12) as completed: Completed Points = 3 and a
#   completion Date of 11/02/2020, plus aligning the Estimation Points (C11,C12)
#   to 3.
# - Moves the last-saved cell selection to F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 1 total points burndown value.
$ws.Range("F2").Value = 17

# Row 10 - Encrypt User Passwords (Estimation Points already 3; completes now)
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "2/11/2020"

# Row 11 - Implement Login Functionality
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 3
$ws.Range("E11").Value = "2/11/2020"

# Row 12 - Add Login UI
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = "2/11/2020"

# Last active selection in the saved file.
$ws.Range("F4").Select() | Out-Null
